$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price cells whose new values look numeric,
# so Excel doesn't auto-convert them to floats (column D is always text,
# matching the rest of the sheet's inline-string cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.718.04'
$ws.Range("D3").Value = '2.647.85'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '537.51'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").Value = '146.53'
$ws.Range("E6").Value = '  +3.73%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").Value = '6.88'
$ws.Range("E9").Value = '  +6.42%  '
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '3.119.67'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").Value = '59.634.60'
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").Value = '21.48'
$ws.Range("E15").Value = '  +4.33%  '
$ws.Range("D16").Value = '2.681.95'
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D19").Value = '340.42'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '66.52'
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("D24").Value = '0.419'
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '7.33'
$ws.Range("E27").Value = '  +1.59%  '
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -3.63%  '
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("D32").Value = '18.89'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").Value = '151.04'
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  +2.90%  '
$ws.Range("D36").Value = '0.840'
$ws.Range("E36").Value = '  +3.07%  '
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("E39").Value = '  +1.77%  '
$ws.Range("D40").Value = '287.04'
$ws.Range("E40").Value = '  +3.58%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '0.606'
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = '0.0540'
$ws.Range("E44").Value = '  +2.88%  '
$ws.Range("D45").Value = '19.30'
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("D46").Value = '0.0948'
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("D48").Value = '1.967.25'
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '18.51'
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '4.58'
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").Value = '111.77'
$ws.Range("E51").Value = '  +0.55%  '
